$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 7, inheriting row 6's formatting (copy/insert mirrors a
# manual "duplicate row" in Excel, which is how B6's style carried onto B7).
$ws.Rows("6").Copy()
$ws.Rows("7").Insert()

# Fill in the new test-case data for "CheckDemeritsTest".
$ws.Range("A7").Value = "CheckDemeritsTest"
$ws.Range("B7").Value = "sa020@mailinator.com"
$ws.Range("C7").Value = "Pa`$`$w0rd"

# Move the active selection, matching the post-edit cursor position.
$ws.Range("F11").Select()
